$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 41 (existing rows 41-46 shift down to 42-47)
$ws.Rows("41:41").Insert()

# Populate the newly inserted row 41 with the new weekly price record
$ws.Cells.Item(41, 1).Value = 1
$ws.Cells.Item(41, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(41, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(41, 4).Value = 44449
$ws.Cells.Item(41, 5).Value = 15
$ws.Cells.Item(41, 6).Value = 100114001
$ws.Cells.Item(41, 7).Value = "Papa"
$ws.Cells.Item(41, 8).Value = "Asterix"
$ws.Cells.Item(41, 9).Value = "1a (cosecha lavada)"
$ws.Cells.Item(41, 10).Value = 1000
$ws.Cells.Item(41, 11).Value = 10000
$ws.Cells.Item(41, 12).Value = 11000
$ws.Cells.Item(41, 13).Value = 10500
$ws.Cells.Item(41, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(41, 15).Value = "Provincia de Melipilla"
$ws.Cells.Item(41, 16).Value = 420
$ws.Cells.Item(41, 17).Value = 25
$ws.Cells.Item(41, 18).Value = "Hortaliza"
